# Read MACRO years directly from config of Excel
# Add a "year" column (D) to the "config" sheet with the MACRO years.

$wb = $excel.ActiveWorkbook

$wsMer = $wb.Worksheets.Item("MERtoPPP")
$wsMer.Activate()
$wsMer.Range("L22").Select() | Out-Null

$ws = $wb.Worksheets.Item("config")
$ws.Activate()

$ws.Range("D1").Value = "year"
$ws.Range("D2").Value = 2020
$ws.Range("D3").Value = 2030
$ws.Range("D4").Value = 2040

$ws.Range("G14").Select() | Out-Null
